$d = $word.ActiveDocument

# The third paragraph holds the _GoBack bookmark (empty text).
$bookmarkPara = $d.Paragraphs.Item(3)
$r = $bookmarkPara.Range
$r.InsertParagraphBefore()
$r.InsertParagraphBefore()

# Now paragraphs 3 and 4 are the two freshly inserted empty paragraphs;
# paragraph 5 is the (original) bookmark paragraph.
$d.Paragraphs.Item(3).Range.Text = "Social network"
$d.Paragraphs.Item(4).Range.Text = "3f295763440ec71da202d2259331dc5ce215df7c"

# Add a trailing empty paragraph after the bookmark paragraph.
$d.Paragraphs.Item(5).Range.InsertParagraphAfter()
